$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark that currently sits on the
#    "Title: Practice All Directions Simultaneously" paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 2) Add a closing period "." right after the ADD/ADHD description
#    sentence ("...I can tell if my concentration has improved"), as its
#    own run (matching the diff, which shows it as a separate <w:r>).
$r = $d.Content.Duplicate
$found = $r.Find.Execute("I can tell if my concentration has improved", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(".")
# Touch a formatting property and revert it so the new text is kept as a
# distinct run instead of being silently merged into the previous one.
$r.Bold = 1
$r.Bold = 0

# 3) Move the "_GoBack" bookmark to wrap the whole "tracking of improved
#    concentration" paragraph (including its paragraph mark), reflecting
#    that this is now the location of the most recent edit.
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*tracking of improved concentration*") {
        $d.Bookmarks.Add("_GoBack", $para.Range)
        break
    }
}
